$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update stats for 2025-12 (row 25)
$ws.Range("B25").Value = 6480
$ws.Range("D25").Value = 6043625
$ws.Range("E25").Value = 932.6581790123457
$ws.Range("F25").Value = 9.998302495331867
$ws.Range("H25").Value = 26.57214494456994
